$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Authors column refresh: re-synced reference metadata (whitespace + minor value refresh) ---
$ws.Range("E2").Value = '[Craig%Brammer%NULL%1,         Stacy%Donohue%NULL%1,         Timothy%Elwell%NULL%1,         Eliza%Fishbein%NULL%1,         D''Ante%Forschino%NULL%1,         Dorothy%Horne%NULL%1,         Buffy%Lloyd-Krejci%NULL%1,         Jessica%Little%NULL%1,         Bistra%Nikiforova%NULL%1,         Elizabeth%Winterbauer%NULL%1]'
$ws.Range("E2").ClearFormats()
$ws.Range("E3").Value = '[Ayesha%Appa%Ayesha.appa@ucsf.edu%1,         Gabriel%Chamie%NULL%1,         Aenor%Sawyer%NULL%1,         Kimberly%Baltzell%NULL%1,         Kathryn%Dippell%NULL%1,         Salu%Ribeiro%NULL%1,         Elias%Duarte%NULL%1,         Joanna%Vinden%NULL%1,         CLIAHUB%Consortium%NULL%1,         Jonathan%Kramer-Feldman%NULL%1,         Shahryar%Rahdari%NULL%1,         Doug%MacIntosh%NULL%1,         Katherine%Nicholson%NULL%1,         Jonathan%Im%NULL%1,         Diane%Havlir%NULL%1,         Bryan%Greenhouse%NULL%1]'
$ws.Range("E3").ClearFormats()
$ws.Range("E4").Value = '[Min Cheol%Chang%NULL%1,         Wan-Seok%Seo%NULL%1,         Donghwi%Park%NULL%1,         Jian%Hur%NULL%0]'
$ws.Range("E4").ClearFormats()
$ws.Range("E5").Value = '[Marci L.%Drees%NULL%1,         Mia A.%Papas%NULL%1,         Terri E.%Corbo%NULL%1,         Kimberly D.%Williams%NULL%1,         Sharon T.%Kurfuerst%NULL%2,         Sharon T.%Kurfuerst%NULL%0]'
$ws.Range("E5").ClearFormats()
$ws.Range("E6").Value = '[Erin F.%Flynn%flynne@email.chop.edu%1,         Elizabeth%Kuhn%NULL%1,         Mohammed%Shaik%NULL%1,         Elizabeth%Tarr%NULL%1,         Nicole%Scattolini%NULL%1,         Allison%Ballantine%NULL%1]'
$ws.Range("E6").ClearFormats()
$ws.Range("E7").Value = '[Scott A.%Goldberg%NULL%1,         Robert A.%Bonacci%NULL%1,         Lucas C.%Carlson%NULL%1,         Charles T.%Pu%NULL%1,         Christine S.%Ritchie%NULL%1]'
$ws.Range("E7").ClearFormats()
$ws.Range("E8").Value = '[Travis%Sanchez%NULL%0,         Sadhu%Panda%NULL%2,         Sadhu%Panda%NULL%0,         Ebrahim%Khajeh%NULL%1,         Alexandra%Halalau%alexandra.halalau@beaumont.edu%2,         Alexandra%Halalau%alexandra.halalau@beaumont.edu%0,         Jeffrey%Ditkoff%NULL%2,         Jeffrey%Ditkoff%NULL%0,         Jessica%Hamilton%NULL%2,         Jessica%Hamilton%NULL%0,         Aryana%Sharrak%NULL%2,         Aryana%Sharrak%NULL%0,         Aimen%Vanood%NULL%2,         Aimen%Vanood%NULL%0,         Amr%Abbas%NULL%2,         Amr%Abbas%NULL%0,         James%Ziadeh%NULL%2,         James%Ziadeh%NULL%0]'
$ws.Range("E8").ClearFormats()
$ws.Range("E10").Value = '[Sang Il%Kim%NULL%1,         Ji Yong%Lee%NULL%2,         Ji Yong%Lee%NULL%0]'
$ws.Range("E10").ClearFormats()
$ws.Range("E12").Value = '[Ki Tae%Kwon%NULL%1,         Jae-Hoon%Ko%NULL%2,         Jae-Hoon%Ko%NULL%0,         Heejun%Shin%NULL%2,         Heejun%Shin%NULL%0,         Minki%Sung%NULL%2,         Minki%Sung%NULL%0,         Jin Yong%Kim%NULL%4,         Jin Yong%Kim%NULL%0]'
$ws.Range("E12").ClearFormats()
$ws.Range("E13").Value = '[Elly%Lee%NULL%1,         Nurul Yaqeen%Mohd Esa%NULL%1,         Tong Ming%Wee%NULL%1,         Chun Ian%Soo%NULL%1]'
$ws.Range("E13").ClearFormats()
$ws.Range("E14").Value = '[Yong Shik%Kwon%NULL%1,         Sun Hyo%Park%NULL%2,         Sun Hyo%Park%NULL%0,         Hyun Jung%Kim%NULL%2,         Hyun Jung%Kim%NULL%0,         Ji Yeon%Lee%NULL%0,         Ji Yeon%Lee%NULL%0,         Mi-ri%Hyun%NULL%2,         Mi-ri%Hyun%NULL%0,         Hyun ah%Kim%NULL%2,         Hyun ah%Kim%NULL%0,         Jae Seok%Park%NULL%0,         Jae Seok%Park%NULL%0]'
$ws.Range("E14").ClearFormats()
$ws.Range("E15").Value = '[Po-Ting%Lin%NULL%1,         Ting-Yuan%Ni%NULL%1,         Tren-Yi%Chen%NULL%1,         Chih-Pei%Su%NULL%1,         Hsiao-Fen%Sun%NULL%1,         Mu-Kuan%Chen%NULL%1,         Chu-Chung%Chou%NULL%1,         Po-Yu%Wang%NULL%1,         Yan-Ren%Lin%H6213.lac@gmail.com%1]'
$ws.Range("E15").ClearFormats()
$ws.Range("E16").Value = '[David A%Lindholm%david.a.lindholm4.mil@mail.mil%1,         John L%Kiley%NULL%2,         John L%Kiley%NULL%0,         Nathan K%Jansen%NULL%1,         Robert T%Hoard%NULL%1,         Matthew R%Bondaryk%NULL%1,         Elizabeth M%Stanley%NULL%1,         Gadiel R%Alvarado%NULL%1,         Ana E%Markelz%NULL%1,         Robert J%Cybulski%NULL%1,         Jason F%Okulicz%NULL%1]'
$ws.Range("E16").ClearFormats()
$ws.Range("E17").Value = '[Charmaine Malenab%Manauis%NULL%1,         Marvin%Loh%marvin.loh@mohh.com.sg%1,         James%Kwan%NULL%1,         John%Chua Mingzhou%NULL%1,         Han Jie%Teo%NULL%1,         David%Teng Kuan Peng%NULL%1,         Shawn%Vasoo Sushilan%NULL%1,         Yee Sin%Leo%NULL%1,         Ang%Hou%NULL%1]'
$ws.Range("E17").ClearFormats()
$ws.Range("E18").Value = '[Rivkees%Scott A.%coreGivesNoEmail%1,        Roberson%Shamarial%coreGivesNoEmail%1]'
$ws.Range("E18").ClearFormats()
$ws.Range("E19").Value = '[Eunhye%Seo%NULL%1,         Eunchan%Mun%NULL%2,         Eunchan%Mun%NULL%0,         Wonsool%Kim%NULL%2,         Wonsool%Kim%NULL%0,         Changhwan%Lee%NULL%2,         Changhwan%Lee%NULL%0]'
$ws.Range("E19").ClearFormats()
$ws.Range("E20").Value = '[Aditya%Shah%shah.aditya@mayo.edu%1,         Douglas%Challener%NULL%1,         Aaron J.%Tande%NULL%1,         Maryam%Mahmood%NULL%1,         John C.%O’Horo%NULL%1,         Elie%Berbari%NULL%1,         Sarah J.%Crane%NULL%1]'
$ws.Range("E20").ClearFormats()
$ws.Range("E21").Value = '[Angie N.%Ton%NULL%1,         Tarang%Jethwa%NULL%1,         Karen%Waters%NULL%2,         Karen%Waters%NULL%0,         Leigh L.%Speicher%NULL%2,         Leigh L.%Speicher%NULL%0,         Dawn%Francis%NULL%2,         Dawn%Francis%NULL%0]'
$ws.Range("E21").ClearFormats()
$ws.Range("E22").Value = '[Kate%Mark%NULL%1,         Katie%Steel%NULL%1,         Janet%Stevenson%NULL%1,         Christine%Evans%NULL%1,         Duncan%McCormick%NULL%1,         Lorna%Willocks%NULL%1,         Alison%McCallum%NULL%1,         Laura%Jones%NULL%1,         Ingolfur%Johannessen%NULL%1,         Kate%Templeton%NULL%1,         Oliver%Koch%NULL%0,         Claire%Mackintosh%NULL%1]'
$ws.Range("E22").ClearFormats()

# --- Row 11 (DOI 10.7861/clinmed.2020-0160): previously unresolved reference, now filled in via CROSSREF ---
$ws.Range("C11").Value = '"Drive-through testing in COVID-19: experience from NHS Lothian"'
$ws.Range("C11").ClearFormats()
$ws.Range("E11").Value = '[Katherine%Hill%xref no email%1,  Robyn%Campbell%xref no email%1,  Callum%Mutch%xref no email%1,  Oliver%Koch%xref no email%1,  Claire%Mackintosh%xref no email%1]'
$ws.Range("E11").ClearFormats()
$ws.Range("F11").Value = '10.7861/clinmed.2020-0160'
$ws.Range("F11").ClearFormats()
$ws.Range("G11").Value = 'CROSSREF'
$ws.Range("G11").ClearFormats()
$ws.Range("H11").Value = '''2023-05-05'
$ws.Range("H11").ClearFormats()
